{"js": "// Replace the division problems in the table with the new values,\n// as described by the diff (25 one-to-one text replacements).\nconst replacements = [\n  [\"313\u00f74=\", \"416\u00f72=\"],\n  [\"443\u00f74=\", \"445\u00f73=\"],\n  [\"249\u00f72=\", \"186\u00f75=\"],\n  [\"496\u00f75=\", \"287\u00f72=\"],\n  [\"162\u00f73=\", \"360\u00f73=\"],\n  [\"881\u00f77=\", \"437\u00f73=\"],\n  [\"704\u00f73=\", \"599\u00f72=\"],\n  [\"804\u00f74=\", \"964\u00f74=\"],\n  [\"794\u00f77=\", \"554\u00f78=\"],\n  [\"323\u00f77=\", \"831\u00f73=\"],\n  [\"861\u00f76=\", \"671\u00f72=\"],\n  [\"506\u00f79=\", \"705\u00f74=\"],\n  [\"525\u00f76=\", \"188\u00f72=\"],\n  [\"249\u00f78=\", \"431\u00f73=\"],\n  [\"129\u00f77=\", \"935\u00f77=\"],\n  [\"698\u00f72=\", \"334\u00f78=\"],\n  [\"629\u00f76=\", \"419\u00f74=\"],\n  [\"955\u00f75=\", \"662\u00f79=\"],\n  [\"917\u00f73=\", \"720\u00f73=\"],\n  [\"736\u00f72=\", \"147\u00f76=\"],\n  [\"198\u00f78=\", \"321\u00f73=\"],\n  [\"259\u00f73=\", \"164\u00f73=\"],\n  [\"386\u00f77=\", \"418\u00f73=\"],\n  [\"127\u00f77=\", \"568\u00f79=\"],\n  [\"992\u00f76=\", \"717\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division problems in the table with the new values,\n# as described by the diff (25 one-to-one text replacements).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"313\u00f74=\", \"416\u00f72=\"),\n    @(\"443\u00f74=\", \"445\u00f73=\"),\n    @(\"249\u00f72=\", \"186\u00f75=\"),\n    @(\"496\u00f75=\", \"287\u00f72=\"),\n    @(\"162\u00f73=\", \"360\u00f73=\"),\n    @(\"881\u00f77=\", \"437\u00f73=\"),\n    @(\"704\u00f73=\", \"599\u00f72=\"),\n    @(\"804\u00f74=\", \"964\u00f74=\"),\n    @(\"794\u00f77=\", \"554\u00f78=\"),\n    @(\"323\u00f77=\", \"831\u00f73=\"),\n    @(\"861\u00f76=\", \"671\u00f72=\"),\n    @(\"506\u00f79=\", \"705\u00f74=\"),\n    @(\"525\u00f76=\", \"188\u00f72=\"),\n    @(\"249\u00f78=\", \"431\u00f73=\"),\n    @(\"129\u00f77=\", \"935\u00f77=\"),\n    @(\"698\u00f72=\", \"334\u00f78=\"),\n    @(\"629\u00f76=\", \"419\u00f74=\"),\n    @(\"955\u00f75=\", \"662\u00f79=\"),\n    @(\"917\u00f73=\", \"720\u00f73=\"),\n    @(\"736\u00f72=\", \"147\u00f76=\"),\n    @(\"198\u00f78=\", \"321\u00f73=\"),\n    @(\"259\u00f73=\", \"164\u00f73=\"),\n    @(\"386\u00f77=\", \"418\u00f73=\"),\n    @(\"127\u00f77=\", \"568\u00f79=\"),\n    @(\"992\u00f76=\", \"717\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)  # wdReplaceAll\n}\n"}
